$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Apply the same cell formatting already used by B38:B45 (wrap-text,
# explicit black font) to the newly-filled-in age values in B46:B65 ---
$ws.Range("B38").Copy() | Out-Null
$ws.Range("B46:B65").PasteSpecial(-4122) | Out-Null

# --- Update the age values for rows 46-65 with the real data ---
# (index 0 -> row 46, index 1 -> row 47, ...)
$ages = @(22, 18, 20, 19, 24, 24, 22, 21, 39, 20, 21, 20, 60, 24, 20, 58, 20, 22, 29, 32)

for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = 46 + $i
    $ws.Cells.Item($row, 2).Value = $ages[$i]
}

# --- Update sheet view: clear the frozen/scrolled topLeftCell and move the
# active selection from D61 to the full row A10:XFD10 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A10:XFD10").Select() | Out-Null

# --- Update workbook view window position ---
$win.Left = 100
$win.Top = 1500
